# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets
# to reflect newly generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 262
$ws1.Range("F5").Value = 3085
$ws1.Range("F6").Value = 2058
$ws1.Range("F9").Value = 1155
$ws1.Range("F11").Value = 908
$ws1.Range("F12").Value = 76

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 262
$ws4.Range("F5").Value = 3085
$ws4.Range("F6").Value = 2058
$ws4.Range("F10").Value = 1155
$ws4.Range("F12").Value = 908
$ws4.Range("F13").Value = 76
